$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.4466192170818505
$ws.Range("C2").Value = 0.07507507507507508
$ws.Range("D2").Value = 0.8928571428571429
$ws.Range("E2").Value = 0.1385041551246537
$ws.Range("F2").Value = 0.2808988764044944
$ws.Range("G2").Value = 0.6292352371732817
$ws.Range("H2").Value = 0.7624398073836276
$ws.Range("I2").Value = 25
$ws.Range("J2").Value = 308
$ws.Range("K2").Value = 226
$ws.Range("L2").Value = 3

# ---- Sheet: Classification Report ----
$ws = $wb.Worksheets.Item("Classification Report")
# Row 2 (label 0)
$ws.Range("B2").Value = 0.9868995633187773
$ws.Range("C2").Value = 0.4232209737827715
$ws.Range("D2").Value = 0.5923984272608126

# Row 3 (label 1)
$ws.Range("B3").Value = 0.07507507507507508
$ws.Range("C3").Value = 0.8928571428571429
$ws.Range("D3").Value = 0.1385041551246537

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.4466192170818505
$ws.Range("C4").Value = 0.4466192170818505
$ws.Range("D4").Value = 0.4466192170818505
$ws.Range("E4").Value = 0.4466192170818505

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.5309873191969262
$ws.Range("C5").Value = 0.6580390583199572
$ws.Range("D5").Value = 0.3654512911927332

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9414705852568135
$ws.Range("C6").Value = 0.4466192170818505
$ws.Range("D6").Value = 0.5697844777593669

# ---- Sheet: Confusion Matrix ----
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 226
$ws.Range("C2").Value = 308
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 25
